# Backlog.xlsx - "Cambios estructurales en ingreso de articulos"
# Se asocia en el detalle de ingreso el detalle de orden de compra
#
# Adds three new backlog rows at the bottom of Hoja1 (rows 92-94),
# each with status "no comenzado" in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$newTasks = @(
    "fix reporte orden de compra calculo de total subtotal e ivas",
    "EN REPORTE PEDIDOS AGREGAR COLUMNA DE OT PARA SABER LO QUE ESTA PROGRAMADO",
    "EN FORMULARIO DE FACTURA PROVEEDOR NO APARECEN LAS OBSERVACIONES GUARDADAS"
)

$startRow = 92
for ($i = 0; $i -lt $newTasks.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTasks[$i]
    $ws.Cells.Item($row, 2).Value = "no comenzado"
}

# Move the view/selection to reflect where the new rows were added.
$ws.Range("C98").Select()
